$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.310.18'
$ws.Range('E2').Value = '  -0.50%  '

$ws.Range('D3').Value = '2.067.14'
$ws.Range('E3').Value = '  +3.28%  '

$ws.Range('E4').Value = '  +0.19%  '

$ws.Range('D5').Value = '''327.04'
$ws.Range('E5').Value = '  +0.80%  '

$ws.Range('E6').Value = '  +0.18%  '

$ws.Range('D7').Value = '''0.5164'
$ws.Range('E7').Value = '  +1.10%  '

$ws.Range('D8').Value = '''0.4306'
$ws.Range('E8').Value = '  +3.31%  '

$ws.Range('D9').Value = '''0.08628'
$ws.Range('E9').Value = '  -1.19%  '

$ws.Range('D10').Value = '''45.69'
$ws.Range('E10').Value = '  +6.22%  '

$ws.Range('D11').Value = '''1.147'
$ws.Range('E11').Value = '  +1.16%  '

$ws.Range('D12').Value = '''23.99'
$ws.Range('E12').Value = '  -2.65%  '

$ws.Range('D13').Value = '2.072.59'
$ws.Range('E13').Value = '  +3.75%  '

$ws.Range('D14').Value = '''6.587'
$ws.Range('E14').Value = '  +0.18%  '

$ws.Range('D15').Value = '''7.600'
$ws.Range('E15').Value = '  +1.94%  '

$ws.Range('D16').Value = '''1.004'
$ws.Range('E16').Value = '  +0.32%  '

$ws.Range('D17').Value = '''94.42'
$ws.Range('E17').Value = '  +0.20%  '

$ws.Range('D18').Value = '''0.00001108'
$ws.Range('E18').Value = '  -0.63%  '

$ws.Range('D19').Value = '''0.06596'
$ws.Range('E19').Value = '  +1.59%  '

$ws.Range('D20').Value = '''18.61'
$ws.Range('E20').Value = '  -1.59%  '

$ws.Range('D22').Value = '''6.179'
$ws.Range('E22').Value = '  -0.15%  '

$ws.Range('D23').Value = '30.353.82'
$ws.Range('E23').Value = '  -0.57%  '

$ws.Range('D24').Value = '''12.17'
$ws.Range('E24').Value = '  +2.78%  '

$ws.Range('D25').Value = '''2.281'
$ws.Range('E25').Value = '  +2.27%  '

$ws.Range('D26').Value = '2.305.38'
$ws.Range('E26').Value = '  +3.35%  '

$ws.Range('E27').Value = '  -1.22%  '

$ws.Range('D28').Value = '''160.10'
$ws.Range('E28').Value = '  -1.95%  '

$ws.Range('D29').Value = '''2.485'
$ws.Range('E29').Value = '  +3.78%  '

$ws.Range('D30').Value = '''129.88'
$ws.Range('E30').Value = '  -1.35%  '

$ws.Range('D31').Value = '''1.165'
$ws.Range('E31').Value = '  +2.54%  '

$ws.Range('E32').Value = '  +1.04%  '

$ws.Range('D33').Value = '''6.022'
$ws.Range('E33').Value = '  -0.75%  '

$ws.Range('D34').Value = '''3.822'
$ws.Range('E34').Value = '  -0.78%  '

$ws.Range('D35').Value = '''1.480'
$ws.Range('E35').Value = '  +10.67%  '

$ws.Range('D36').Value = '''0.02536'
$ws.Range('E36').Value = '  +0.52%  '

$ws.Range('D37').Value = '''9.461'
$ws.Range('E37').Value = '  +4.75%  '

$ws.Range('D38').Value = '''5.387'
$ws.Range('E38').Value = '  -1.14%  '

$ws.Range('D39').Value = '''0.06554'
$ws.Range('E39').Value = '  -0.76%  '

$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').Value = '''12.38'
$ws.Range('E40').Value = '  -0.37%  '

$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '''0.2217'
$ws.Range('E41').Value = '  +0.80%  '

$ws.Range('D42').Value = '''0.6616'
$ws.Range('E42').Value = '  -0.07%  '

$ws.Range('D43').Value = '''1.226'
$ws.Range('E43').Value = '  -0.69%  '

$ws.Range('E44').Value = '  +0.12%  '

$ws.Range('D45').Value = '''13.86'
$ws.Range('E45').Value = '  +1.84%  '

$ws.Range('D46').Value = '''0.6236'
$ws.Range('E46').Value = '  +1.13%  '

$ws.Range('D47').Value = '''2.174'
$ws.Range('E47').Value = '  -1.35%  '

$ws.Range('D48').Value = '''3.603'
$ws.Range('E48').Value = '  -1.80%  '

$ws.Range('D49').Value = '''1.227'

$ws.Range('D50').Value = '''80.98'
$ws.Range('E50').Value = '  +0.69%  '

$ws.Range('D51').Value = '''1.171'
$ws.Range('E51').Value = '  +5.70%  '
